# Update FuelPrices at 2025-04-14 03:55
#
# The sheet had its last data row at row 30 (A30:C30), where C30 used a
# "date only" number format (YYYY-MM-DD) that differed from the
# "date + time" number format (YYYY-MM-DD HH:MM:SS) used by the rest of
# column C (e.g. C29). A new data point has arrived, so:
#   1. Row 30 becomes a regular data row: C30 gets the standard
#      datetime number format used throughout column C.
#   2. A new row 31 is appended with the latest values, and it is C31
#      that now carries the special "date only" number format that used
#      to sit on C30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Give C30 the same number format as the rest of column C (C29).
$ws.Range("C30").NumberFormat = $ws.Range("C29").NumberFormat()

# 2) Append the new row 31, carrying forward the latest observed values.
$ws.Range("A31").Value = $ws.Range("A30").Value()
$ws.Range("B31").Value = $ws.Range("B30").Value()
$ws.Range("C31").Value = $ws.Range("C30").Value2()

# The new last row takes over the "date only" formatting that used to be
# on C30.
$ws.Range("C31").NumberFormat = "YYYY-MM-DD"
